# "delete employee date in excel"
#
# The "Employment date" column (with its "DD/MM/YYYY" placeholder data)
# is removed entirely from the import template. Deleting the whole column
# shifts every column to its right one position to the left, which also
# folds the now-unused "Employment date" / "DD/MM/YYYY" shared strings out
# of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRow = 1
$lastCol = $ws.UsedRange.Columns.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item($headerRow, $c)
    if ($cell.Value2 -eq "Employment date") {
        $cell.EntireColumn.Delete()
        break
    }
}
